# Update "县文化馆" (County Culture Center) yearly statistics sheet.
# The 2008 and 2009 rows are removed, the 2010-2020 rows shift up
# accordingly, and a new 2021 row is appended at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two oldest years (2008 and 2009), which occupied rows 2 and 3.
# Everything below (2010..2020) shifts up by two rows automatically.
$ws.Rows("2:3").Delete()

# After the delete, the last data row (2020) is row 12, and the new
# row for 2021 goes to row 13. Copy the formatting of the year label
# cell (A12, which carries the bold/bordered/centered style) onto the
# new label cell A13 before filling in the 2021 data.
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)

$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 14996
$ws.Range("C13").Value = 83137
$ws.Range("D13").Value = 21478
$ws.Range("E13").Value = 3662
$ws.Range("F13").Value = 48511
$ws.Range("G13").Value = 336
$ws.Range("H13").Value = 586.29
$ws.Range("I13").Value = 12342.33
$ws.Range("J13").Value = 231.84
$ws.Range("K13").Value = 2150.35
$ws.Range("L13").Value = 22514
$ws.Range("M13").Value = 358143.6
$ws.Range("N13").Value = 35.07961
$ws.Range("O13").Value = 1585
$ws.Range("P13").Value = 12756
$ws.Range("Q13").Value = 123418
$ws.Range("R13").Value = 42126
